$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.103.71"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.786.64"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  -3.40%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "2.044.39"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "1.787.95"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "34.048.80"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.619"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").Value = "1.449.79"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.646"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.40%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0192"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.918"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0507"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "0.0₆0136"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").Value = "1.945.66"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +0.20%  "
